$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.24 = 12410.2 pesos`n✅ 12410.2 pesos = 3.24 = 977.54 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 309.09
$ws2.Range("O10").Value = 3835.87
$ws2.Range("N12").Value = 3833.99
$ws2.Range("O12").Value = 302
